$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.078.06'
$ws.Range("E2").Value = '  +5.53%  '

$ws.Range("D3").Value = '1.917.66'
$ws.Range("E3").Value = '  +2.52%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.61%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '329.99'
$ws.Range("D5").Style = "Normal"

$ws.Range("E6").Value = '  -0.56%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5217'
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4087'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.76%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08495'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.65%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.00'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.94%  '

$ws.Range("E11").Value = '  +1.96%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.56'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +10.73%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.424'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.39%  '

$ws.Range("D14").Value = '1.926.75'
$ws.Range("E14").Value = '  +2.85%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.425'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.02%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.002'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.60%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '95.45'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.82%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001112'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.84%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06719'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.33'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.34%  '

$ws.Range("E21").Value = '  -0.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.010'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.59%  '

$ws.Range("D23").Value = '30.096.98'
$ws.Range("E23").Value = '  +5.48%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.32'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.80%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.220'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.60%  '

$ws.Range("D26").Value = '2.153.64'
$ws.Range("E26").Value = '  +3.27%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.55'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.16%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.09'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.28%  '

$ws.Range("E29").Value = '  +1.56%  '

$ws.Range("E30").Value = '  +2.34%  '

$ws.Range("E31").Value = '  +3.50%  '

$ws.Range("E32").Value = '  +1.47%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.082'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.99%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.636'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.44%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02488'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.21%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06634'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.87%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2212'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.20%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.231'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.40%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.189'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.22%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.906'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.08%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6534'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.65%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.64'
$ws.Range("D42").Style = "Normal"

$ws.Range("E43").Value = '  +0.72%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6166'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.66%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.31'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.61%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.768'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.30%  '

$ws.Range("E47").Value = '  +4.00%  '

$ws.Range("E48").Value = '  +2.66%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.63'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.07%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.164'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +12.58%  '

$ws.Range("E51").Value = '  +4.54%  '
